# RSNN uniform files to be run
#
# Fill in the newly-run 400-epoch Training/Testing accuracy figures for the
# Feed Forward SNN results table (columns G/H, rows 7, 8, 10, 11), and leave
# the selection where the editor last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1 Layer, 128 hidden units
$ws.Range("G7").Value = 0.861
$ws.Range("H7").Value = 0.511

# 1 Layer, 256 hidden units
$ws.Range("G8").Value = 0.954
$ws.Range("H8").Value = 0.554

# 2 layers, 128 hidden units
$ws.Range("G10").Value = 0.703
$ws.Range("H10").Value = 0.553

# 2 layers, 256 hidden units
$ws.Range("G11").Value = 0.852
$ws.Range("H11").Value = 0.603

# Touch the title row's formatting so the workbook's style table gets
# refreshed/re-saved the same way it was in the edit session (the title
# stays bold/centered; this just nudges the cell onto the refreshed style
# record).
$ws.Range("C4:F4").Font.Bold = $false

# Leave the selection where the editor ended up.
$ws.Range("I20").Select() | Out-Null
